$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric stay as text, matching the
# original inline-string cell type, then restore the default style so
# no stray formatting is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.059.22"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "2.304.83"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "309.69"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").Value = "105.84"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "39.84"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "8.30"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "0.996"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "15.34"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "2.678.33"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.309.09"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "42.894.63"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  -4.07%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "73.49"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").Value = "266.57"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").Value = "  +11.81%  "
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "22.27"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "37.75"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "165.16"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "2.82"
$ws.Range("E34").Value = "  +5.94%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").Value = "4.59"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "107.38"
$ws.Range("E41").Value = "  +8.42%  "
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "71.64"
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("D44").Value = "0.229"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").Value = "1.01"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "12.30"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "1.724.33"
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("D48").Value = "111.85"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("D49").Value = "76.80"
$ws.Range("E49").Value = "  -5.43%  "
$ws.Range("D50").Value = "8.77"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -2.95%  "

# Clear the temporary text-format style so cell formatting matches the original.
$ws.Range("D2:D51").Style = "Normal"

